$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header info changes
$ws.Range("B2").Value = "Agile Co-Development Services"
$ws.Range("B3").Value = "GVT000ABC1234"
$ws.Range("B7").Value = "asdasdad"
$ws.Range("B8").Value = "asdasdConsulting"

# Move "At Work" marks (column C) to "Annual Leave" marks (column G) for specific days
$rows = @(31, 32, 33, 34, 37, 38, 41)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 7).Value = 1
}

# Totals row 44
$ws.Range("C44").Value = 12
$ws.Range("G44").Value = 7

# Footer changes
$ws.Range("B48").Value = "12 - Feb - 2025"
$ws.Range("B50").Value = "John Doe"
